# "fix issue with mapping" - Sheet1 had its Question/Answer/Intent/Citations
# mapping straightened out: the OSI question now maps to the out-of-domain
# canned answer, citations were filled in for the in-domain rows, and a new
# out-of-domain example row (joke about Firewall) was appended.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (use_app_gw / Tell me about Azure Application Gateway): answer now
# carries a [doc0] citation marker, and the Citations column is populated.
$ws.Range('C2').Value = 'Azure Application Gateway is a web traffic load balancer that enables you to manage traffic to your web application. It operates at the application layer (OSI layer 7) and routes traffic based on the content of the request, such as the URL or HTTP header. It features a web application firewall and intelligent layer 7 routing. [doc0] '
$ws.Range('E2').Value = '[doc0],What is Azure Firewall?,https://learn.microsoft.com/en-us/azure/firewall/overview'

# Row 3 (use_app_gw / What is OSI?): previously answered in-domain with the
# OSI explanation; now correctly mapped to the out-of-domain canned reply.
$ws.Range('C3').Value = 'Sorry, I am a chat bot that can answer questions on Azure.'
$ws.Range('D3').Value = 'OUT_OF_DOMAIN'

# Row 4 (use_app_gw / private connection question): answer now cites
# [doc0], and the Citations column points at the Private Link doc.
$ws.Range('C4').Value = 'Yes [doc0]'
$ws.Range('E4').Value = '[doo0],Application Gateway Private Link,https://learn.microsoft.com/en-us/azure/application-gateway/private-link'

# Row 5 (use_waf / firewall solutions question) shifts down from the old
# row 5 content and gains the same Citations value as row 2.
$ws.Range('A5').Value = 'use_waf'
$ws.Range('B5').Value = 'What firewall solution exist in Azure?'
$ws.Range('C5').Value = 'Azure Firewall is a cloud-native and intelligent network firewall security service that provides the best of breed threat protection for your cloud workloads running in Azure. It''s a fully stateful firewall as a service with built-in high availability and unrestricted cloud scalability. It provides both east-west and north-south traffic inspection. [doc0]'
$ws.Range('D5').Value = 'IN_DOMAIN'
$ws.Range('E5').Value = '[doc0],What is Azure Firewall?,https://learn.microsoft.com/en-us/azure/firewall/overview'

# New row 6: an out-of-domain example (joke about Firewall -> canned reply).
$ws.Range('A6').Value = 'out_of_domain'
$ws.Range('B6').Value = 'Tell me a joke about Firewall'
$ws.Range('C6').Value = 'Sorry, I am a chat bot that can answer questions on Azure.'
$ws.Range('D6').Value = 'OUT_OF_DOMAIN'

# Selection moved as a side effect of editing further down the sheet.
$ws.Range('C11').Select() | Out-Null
